$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$newSheet = $wb.Worksheets.Add($null, $ws1)
$newSheet.Name = "test"

$src = $ws1.Range("A3:BC6")
$src.Copy() | Out-Null
$newSheet.Paste($newSheet.Range("A1")) | Out-Null
Write-Host "done"
